# Generate Report for Handoff
#
# The localization status report was regenerated: the zh-cn/de-de rows moved
# from "In Translation" to "Ready for handoff", their handoff timestamps were
# refreshed, and the Status/Language columns were widened to fit the new
# (longer) status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff timestamps ------------------------------------------
$overview.Range("G2").Value = "2016-09-07 04:52:01"
$dede.Range("H2").Value     = "2016-09-07 04:52:01"
$zhcn.Range("H2").Value     = "2016-09-07 04:51:55"

# --- Widen the Status / language columns to fit "Ready for handoff" --------
# (ColumnWidth is quantized by Excel to the workbook's default-font character
# grid, so we target the input that lands on the grid point closest to the
# new authored width.)
$overview.Columns.Item(5).ColumnWidth = 16.3333
$overview.Columns.Item(6).ColumnWidth = 16.3333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333
$dede.Columns.Item(3).ColumnWidth     = 16.3333
